# Insert a new data row at row 35 of the active sheet, shifting the
# existing rows 35-62 down to 36-63 (this matches the diff: every row
# from the old 35..62 reappears, unchanged, one row further down, and a
# brand-new record is inserted at row 35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 35:62 down to 36:63, leaving a blank row 35 behind.
$ws.Rows(35).Insert()

# Populate the newly inserted row 35 with the new record's data.
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 44484
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100108
$ws.Range("H35").Value = "Tropicales y subtropicales"
$ws.Range("I35").Value = 100108002
$ws.Range("J35").Value = "Mango"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 60
$ws.Range("N35").Value = 7500
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 7750
$ws.Range("Q35").Value = "$/bandeja 4 kilos"
$ws.Range("R35").Value = "Perú"
$ws.Range("S35").Value = 1938
$ws.Range("T35").Value = 4
